# Bot5GUI -v1.0 "intermedio" update
# Deletes the "Tipo de cuenta" / "Validacion" configuration rows
# (rows 9 and 10) from the "parametrosInicio" sheet, shifting the
# remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")

# Select rows 9 and 10 (Tipo de cuenta / Validacion) and delete them
# entirely, shifting rows below upward.
$rows = $ws.Range("A9:B10").EntireRow
$rows.Delete()

# Row 8 ("Layer" / "/PDDPROD") is now the last row of its bordered
# box, so it picks up the thicker "closing" bottom border that Excel
# applies to the final row of a box (matching the box previously
# closed out by the now-deleted rows).
$box = $ws.Range("A8:B8")
$box.Borders.Item(9).LineStyle = 1      # xlEdgeBottom
$box.Borders.Item(9).Weight = -4138     # xlMedium

# Update the active selection/view on the sheet as recorded after edit.
$ws.Activate()
$ws.Range("D8").Select()

$wb.Save()
